# Finished combining the survival parts
# Populate the AC/MI (weighted / unweighted) HR + CI figures for the new
# "Cape. Vs none" / "Other vs. none" / "Cape vs Other" rows (46-48) on the
# "cape cox everything" worksheet, and update the existing "MI weights"
# column (M/N) values to the newly-rounded figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cape cox everything")
$ws.Activate()

# --- numeric (HR) columns: D / G / J / M for rows 46-48 -------------------
$ws.Range("D46").Value = 0.687
$ws.Range("G46").Value = 0.603
$ws.Range("J46").Value = 0.752
$ws.Range("M46").Value = 0.701

$ws.Range("D47").Value = 0.521
$ws.Range("G47").Value = 0.452
$ws.Range("J47").Value = 0.579
$ws.Range("M47").Value = 0.532

$ws.Range("D48").Value = 1.318
$ws.Range("G48").Value = 1.334
$ws.Range("M48").Value = 1.317

# --- confidence-interval text columns: E / H / K / N ----------------------
# (written in the same order the original author entered them so new shared
# strings line up the same way)
$ws.Range("N48").Value = "(1.100,1.579)"
$ws.Range("H48").Value = "(1.109,1.604)"
$ws.Range("H46").Value = "(0.443,0.820)"
$ws.Range("H47").Value = "(0.340,0.602)"
$ws.Range("E46").Value = "(0.530,0.891)"
$ws.Range("E47").Value = "(0.416,0.653)"
$ws.Range("E48").Value = "(1.078,1.612)"
$ws.Range("K46").Value = "(0.595,0.952)"
$ws.Range("K47").Value = "(0.474,0.707)"

# J48 keeps a trailing zero ("1.300"), so it has to be stored as text
# (format it as Text first so Excel doesn't coerce it back to the number 1.3).
$ws.Range("J48").NumberFormat = "@"
$ws.Range("J48").Value = "1.300"

$ws.Range("K48").Value = "(1.076,1.570)"

# --- view state: scroll position + active selection ------------------------
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("N43").Select()
